$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Row 1: turn the old (duplicated data) row into a proper header row ----
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "capacity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "register_date"
$ws.Cells.Item(1,6).Value = "register_reason"
$ws.Cells.Item(1,7).Value = "acquire_value"
$ws.Cells.Item(1,8).Value = "property_category"
$ws.Cells.Item(1,9).Value = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = "legislator_id"
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = "index"

# Give the newly added header cells (H1:N1) the same look (bold / border) as the rest of row 1
$ws.Range("B1:G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Row 2: keep the existing data values, fix B2/E2 and append the metadata columns ----
$ws.Cells.Item(2,2).Value = "HondaCRV"
$ws.Cells.Item(2,5).Value = "96年07月11曰"

# Give the newly added data cells (H2:N2) the same look as the rest of row 2 first,
# then fill in their values (so the later NumberFormat tweak below is not overwritten)
$ws.Range("B2:G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(2,8).Value = "land"
$ws.Cells.Item(2,9).Value = "normal"
$ws.Cells.Item(2,10).NumberFormat = "@"
$ws.Cells.Item(2,10).Value = "2012-04-24"
$ws.Cells.Item(2,11).Value = "蔡其昌"
$ws.Cells.Item(2,12).Value = 1377
$ws.Cells.Item(2,13).Value = "tmp61ee1"
$ws.Cells.Item(2,14).Value = 39
